# Apply crypto price/volume updates from the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.293.24'
$ws.Range("E2").Value = '  +2.13%  '

$ws.Range("D3").Value = '1.801.54'
$ws.Range("E3").Value = '  +3.58%  '

$ws.Range("E4").Value = '  -0.39%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.68'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.36%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4586'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +19.78%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3803'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +13.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.22'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.153'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.10%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07584'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.51'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.16%  '

$ws.Range("E13").Value = '  -0.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.353'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.21%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.571'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.52%  '

$ws.Range("D16").Value = '1.804.91'
$ws.Range("E16").Value = '  +3.36%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001093'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.37%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06727'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '81.52'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.96%  '

$ws.Range("E20").Value = '  -0.56%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.46'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.419'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.52%  '

$ws.Range("D23").Value = '28.291.02'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.90'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.428'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.69'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.96'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.366'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +5.18%  '

$ws.Range("D29").Value = '2.011.55'
$ws.Range("E29").Value = '  +3.54%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.99'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.82%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.248'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.028'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09501'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +10.45%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.860'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.40%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2309'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.11'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.76%  '

$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02357'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.51%  '

$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.279'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06351'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.32%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6628'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.27%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.238'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.378'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.77%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.486'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.94%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.28'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.54%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.874'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.82%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.6129'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.21'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.039'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07167'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.52%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.177'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.13%  '
